$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3008.3333
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2685
$ws.Range("H79").Value = 3008.3333
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1908
$ws.Range("H86").Value = 4580.8
$ws.Range("I86").Value = 4693.4165
$ws.Range("J86").Value = 4476.846
$ws.Range("K86").Value = 4693.4165
$ws.Range("L86").Value = 4476.846
$ws.Range("M86").Value = -3570.4165
$ws.Range("N86").Value = -6722.846
$ws.Range("H89").Value = 4580.8
$ws.Range("I89").Value = 4693.4165
$ws.Range("J89").Value = 4476.846
$ws.Range("K89").Value = 23467.0825
$ws.Range("L89").Value = 22384.23
$ws.Range("M89").Value = -17851.0825
$ws.Range("N89").Value = -33616.23
$ws.Range("H137").Value = 1247.5667
$ws.Range("I137").Value = 1090.6296
$ws.Range("J137").Value = 2660
$ws.Range("K137").Value = 3271.8888
$ws.Range("L137").Value = 7980
$ws.Range("M137").Value = -721.8887999999997
$ws.Range("N137").Value = -13080
$ws.Range("H138").Value = 2411484.5
$ws.Range("J138").Value = 11116392
$ws.Range("L138").Value = 33349176
$ws.Range("N138").Value = -33359456
$ws.Range("H141").Value = 595.11536
$ws.Range("I141").Value = 595.11536
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1785.34608
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3394.65392
$ws.Range("N141").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 19250
$ws.Range("J59").Value = 19250
$ws.Range("L59").Value = 19250
$ws.Range("N59").Value = -20858
$ws.Range("H61").Value = 1192.6586
$ws.Range("I61").Value = 1141.3077
$ws.Range("J61").Value = 2194
$ws.Range("K61").Value = 1141.3077
$ws.Range("L61").Value = 2194
$ws.Range("M61").Value = -929.3077000000001
$ws.Range("N61").Value = -2618
$ws.Range("H74").Value = 24287.701
$ws.Range("I74").Value = 32036.219
$ws.Range("J74").Value = 7757.533
$ws.Range("K74").Value = 32036.219
$ws.Range("L74").Value = 7757.533
$ws.Range("M74").Value = -31162.219
$ws.Range("N74").Value = -9505.532999999999
$ws.Range("H77").Value = 24287.701
$ws.Range("I77").Value = 32036.219
$ws.Range("J77").Value = 7757.533
$ws.Range("K77").Value = 160181.095
$ws.Range("L77").Value = 38787.665
$ws.Range("M77").Value = -155813.095
$ws.Range("N77").Value = -47523.665
$ws.Range("H88").Value = 2500.6
$ws.Range("J88").Value = 2766.6667
$ws.Range("L88").Value = 2766.6667
$ws.Range("N88").Value = -3578.6667
$ws.Range("H91").Value = 2500.6
$ws.Range("J91").Value = 2766.6667
$ws.Range("L91").Value = 2766.6667
$ws.Range("N91").Value = -5574.6667
$ws.Range("H136").Value = 1192.6586
$ws.Range("I136").Value = 1141.3077
$ws.Range("J136").Value = 2194
$ws.Range("K136").Value = 3423.9231
$ws.Range("L136").Value = 6582
$ws.Range("M136").Value = -873.9231
$ws.Range("N136").Value = -11682

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29169.5
$ws.Range("J55").Value = 29169.5
$ws.Range("L55").Value = 29169.5
$ws.Range("N55").Value = -29715.5
$ws.Range("H86").Value = 1739.8636
$ws.Range("I86").Value = 1604.6666
$ws.Range("J86").Value = 2029.5714
$ws.Range("K86").Value = 1604.6666
$ws.Range("L86").Value = 2029.5714
$ws.Range("M86").Value = -481.6666
$ws.Range("N86").Value = -4275.5714
$ws.Range("H89").Value = 1739.8636
$ws.Range("I89").Value = 1604.6666
$ws.Range("J89").Value = 2029.5714
$ws.Range("K89").Value = 8023.333000000001
$ws.Range("L89").Value = 10147.857
$ws.Range("M89").Value = -2407.333000000001
$ws.Range("N89").Value = -21379.857
$ws.Range("H134").Value = 730632.2
$ws.Range("I134").Value = 1252347.4
$ws.Range("J134").Value = 4767.6523
$ws.Range("K134").Value = 3757042.2
$ws.Range("L134").Value = 14302.9569
$ws.Range("M134").Value = -3754507.2
$ws.Range("N134").Value = -19372.9569

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 429.32
$ws.Range("I22").Value = 257.14285
$ws.Range("J22").Value = 648.4545000000001
$ws.Range("K22").Value = 257.14285
$ws.Range("L22").Value = 648.4545000000001
$ws.Range("M22").Value = 92.85714999999999
$ws.Range("N22").Value = -1348.4545
$ws.Range("H25").Value = 5108.3335
$ws.Range("I25").Value = 1255.5555
$ws.Range("J25").Value = 16666.666
$ws.Range("K25").Value = 1255.5555
$ws.Range("L25").Value = 16666.666
$ws.Range("M25").Value = -1081.5555
$ws.Range("N25").Value = -17014.666
$ws.Range("H62").Value = 4233.3335
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4233.3335
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4233.3335
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5481.3335
$ws.Range("H65").Value = 4233.3335
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4233.3335
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 21166.6675
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -27406.6675
$ws.Range("H132").Value = 700359.4
$ws.Range("I132").Value = 1497.9111
$ws.Range("J132").Value = 4631455
$ws.Range("K132").Value = 4493.7333
$ws.Range("L132").Value = 13894365
$ws.Range("M132").Value = -1963.7333
$ws.Range("N132").Value = -13899425

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4943.6665
$ws.Range("I63").Value = 2385.3333
$ws.Range("J63").Value = 7502
$ws.Range("K63").Value = 7155.999899999999
$ws.Range("L63").Value = 22506
$ws.Range("M63").Value = -6406.999899999999
$ws.Range("N63").Value = -24004
$ws.Range("H66").Value = 4943.6665
$ws.Range("I66").Value = 2385.3333
$ws.Range("J66").Value = 7502
$ws.Range("K66").Value = 21467.9997
$ws.Range("L66").Value = 67518
$ws.Range("M66").Value = -17723.9997
$ws.Range("N66").Value = -75006
$ws.Range("H131").Value = 867.89
$ws.Range("I131").Value = 535.4545000000001
$ws.Range("J131").Value = 908.97754
$ws.Range("K131").Value = 1606.3635
$ws.Range("L131").Value = 2726.93262
$ws.Range("M131").Value = 3433.6365
$ws.Range("N131").Value = -12806.93262

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 15361.385
$ws.Range("J57").Value = 16516.5
$ws.Range("L57").Value = 16516.5
$ws.Range("N57").Value = -18156.5
$ws.Range("H63").Value = 23616.5
$ws.Range("J63").Value = 24040
$ws.Range("L63").Value = 24040
$ws.Range("N63").Value = -25412
$ws.Range("H66").Value = 23616.5
$ws.Range("J66").Value = 24040
$ws.Range("L66").Value = 72120
$ws.Range("N66").Value = -78984
$ws.Range("H135").Value = 31982.223
$ws.Range("J135").Value = 31982.223
$ws.Range("L135").Value = 31982.223
$ws.Range("N135").Value = -42122.223

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7276.2856
$ws.Range("I68").Value = 11750.2
$ws.Range("J68").Value = 3209.0908
$ws.Range("K68").Value = 11750.2
$ws.Range("L68").Value = 3209.0908
$ws.Range("M68").Value = -11001.2
$ws.Range("N68").Value = -4707.0908
$ws.Range("H71").Value = 7276.2856
$ws.Range("I71").Value = 11750.2
$ws.Range("J71").Value = 3209.0908
$ws.Range("K71").Value = 58751
$ws.Range("L71").Value = 16045.454
$ws.Range("M71").Value = -55007
$ws.Range("N71").Value = -23533.454
$ws.Range("H136").Value = 1420.625
$ws.Range("I136").Value = 1185.909
$ws.Range("J136").Value = 4002.5
$ws.Range("K136").Value = 3557.727
$ws.Range("L136").Value = 12007.5
$ws.Range("M136").Value = -1007.727
$ws.Range("N136").Value = -17107.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27333.334
$ws.Range("J64").Value = 27333.334
$ws.Range("L64").Value = 27333.334
$ws.Range("N64").Value = -27829.334
$ws.Range("H67").Value = 27333.334
$ws.Range("J67").Value = 27333.334
$ws.Range("L67").Value = 27333.334
$ws.Range("N67").Value = -29049.334
$ws.Range("H136").Value = 1421
$ws.Range("I136").Value = 976.125
$ws.Range("K136").Value = 2928.375
$ws.Range("M136").Value = -378.375
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
